$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2: "Sara" -> "Bbzinha"
$ws.Range("A2").Value = "Bbzinha"

# Update B2: phone number, stored as General number instead of Text
$ws.Range("B2").NumberFormat = "general"
$ws.Range("B2").Value = 558586441988

# Remove the stray C9 cell so the used range / dimension shrinks back to A1:B2
$ws.Range("C9").Clear()
